# Add two new test steps (rows 15 and 16) to the TC32 store-location
# testcase sheet, continuing the existing step pattern of
# CLICK / <object> / xpath used throughout the sheet, for the
# "adding 20 items into cart" flow: open the My Account section, then Logout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply the same box-border formatting used by the existing data rows
# (e.g. row 14) to the two new rows before filling in their values.
$ws.Range("A15:E16").Borders.LineStyle = 1

# Row 15: CLICK on MyaccountSection (located via xpath)
$ws.Cells.Item(15, 2).Value2 = "CLICK"
$ws.Cells.Item(15, 3).Value2 = "MyaccountSection"
$ws.Cells.Item(15, 4).Value2 = "xpath"

# Row 16: CLICK on Logout (located via xpath)
$ws.Cells.Item(16, 2).Value2 = "CLICK"
$ws.Cells.Item(16, 3).Value2 = "Logout"
$ws.Cells.Item(16, 4).Value2 = "xpath"

# Match the selection Excel leaves behind after editing the new rows
[void]$ws.Range("A15:XFD16").Select()
